$wb = $excel.ActiveWorkbook

# Add the new worksheet and rename it to "k_p"
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "k_p"

# Move it to the end, after the last existing sheet (soil_spectra)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# Populate header row
$ws = $wb.Worksheets("k_p")
$ws.Range("A1").Value = "P"
$ws.Range("B1").Value = "k_p"

# Populate data rows
$data = @(
    @(2500, 1.8),
    @(2000, 1.7),
    @(1500, 1.5),
    @(1000, 1.3),
    @(800, 1.25),
    @(500, 1),
    @(250, 0.75),
    @(200, 0.7),
    @(100, 0.5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
}

# Select cell B11 on the new sheet (matches the saved selection state)
$ws.Range("B11").Select()

# Make the new sheet the active / selected tab
$wb.Worksheets("k_p").Activate()
